$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spelling/naming inconsistency in the merge-field placeholders:
# "ARCHERRESULT" -> "RESULTS" (matching the __Results__ defined name).
$ws.Range("C6").Value = "<#RESULTS.COMPETITIONNAME>"
$ws.Range("D6").Value = "<#RESULTS.COMPETITIONRESULT>"

# Update the active selection to match the author's saved view.
$ws.Range("C6").Select()
